$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.960.12"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.750.64"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'336.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.3852"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'0.3404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("D10").Value = "'1.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "'0.07224"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").Value = "'22.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'6.167"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "'7.118"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "1.752.64"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "'0.00001061"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "'0.06613"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'79.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'16.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("D22").Value = "'6.186"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").Value = "27.982.56"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "'11.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").Value = "'2.391"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "'153.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").Value = "'2.299"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.17%  "
$ws.Range("D29").Value = "1.951.95"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "'1.264"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.58%  "
$ws.Range("D31").Value = "'131.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").Value = "'4.030"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").Value = "'5.845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").Value = "'0.08828"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "'12.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.07%  "
$ws.Range("D36").Value = "'1.540"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("D37").Value = "'0.6554"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "'0.02283"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.53%  "
$ws.Range("D39").Value = "'5.136"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("D40").Value = "'0.06142"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "'0.2096"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "'1.207"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").Value = "'8.014"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'13.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").Value = "'3.839"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'0.6036"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "'126.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").Value = "'2.004"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").Value = "'1.171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").Value = "'1.107"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.19%  "
